$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Section: Inventory (rows 27-32)
# ---------------------------------------------------------------------------
$ws.Range("A27:F27").Merge()
$ws.Range("A27:F27").Style = "Good"
$ws.Range("A27").Value = "Inventory"

$ws.Range("A28").Value = "Attribute Name"
$ws.Range("B28").Value = "Required"
$ws.Range("C28").Value = "Type"
$ws.Range("D28").Value = "Max Field Size"
$ws.Range("E28").Value = "Primary Key"
$ws.Range("F28").Value = "Extra Info"

$ws.Range("A29").Value = "Item_ID"
$ws.Range("B29").Value = "Yes"
$ws.Range("C29").Value = "Char"
$ws.Range("D29").Value = 6
$ws.Range("E29").Value = "Yes"

$ws.Range("A30").Value = "Item_Description"
$ws.Range("B30").Value = "Yes"
$ws.Range("C30").Value = "Char"
$ws.Range("D30").Value = 30
$ws.Range("E30").Value = "No"

$ws.Range("A31").Value = "Item_Price"
$ws.Range("B31").Value = "Yes"
$ws.Range("C31").Value = "Numerical"
$ws.Range("D31").Value = 6.2
$ws.Range("D31").HorizontalAlignment = -4152
$ws.Range("E31").Value = "no"

$ws.Range("A32").Value = "OnHand"
$ws.Range("B32").Value = "Yes"
$ws.Range("C32").Value = "numerical"
$ws.Range("D32").Value = 6
$ws.Range("E32").Value = "no"

$loInventory = $ws.ListObjects.Add(1, $ws.Range("A28:F36"), 0, 1)
$loInventory.Name = "Table134"
$loInventory.TableStyle = "TableStyleMedium14"

# ---------------------------------------------------------------------------
# Section: Invoices (rows 40-46)
# ---------------------------------------------------------------------------
$ws.Range("A40:F40").Merge()
$ws.Range("A40:F40").Style = "Good"
$ws.Range("A40").Value = "Invoices"

$ws.Range("A41").Value = "Attribute Name"
$ws.Range("B41").Value = "Required"
$ws.Range("C41").Value = "Type"
$ws.Range("D41").Value = "Max Field Size"
$ws.Range("E41").Value = "Primary Key"
$ws.Range("F41").Value = "Extra Info"

$ws.Range("A42").Value = "Invo_Num"
$ws.Range("B42").Value = "Yes"
$ws.Range("C42").Value = "Num"
$ws.Range("D42").Value = 8
$ws.Range("E42").Value = "Yes"

$ws.Range("A43").Value = "Date"
$ws.Range("B43").Value = "Yes"
$ws.Range("C43").Value = "Num"
$ws.Range("D43").Value = 8
$ws.Range("E43").Value = "No"

$ws.Range("A44").Value = "Item_ID"
$ws.Range("B44").Value = "yes"
$ws.Range("C44").Value = "Num"
$ws.Range("D44").Value = 6
$ws.Range("D44").HorizontalAlignment = -4152
$ws.Range("E44").Value = "No"
$ws.Range("F44").Value = "Foreign key"

$ws.Range("A45").Value = "Amount"
$ws.Range("B45").Value = "Yes"
$ws.Range("C45").Value = "Num"
$ws.Range("D45").Value = 7
$ws.Range("E45").Value = "no"

$ws.Range("A46").Value = "Total_Sale"
$ws.Range("B46").Value = "Yes"
$ws.Range("C46").Value = "num"
$ws.Range("D46").Value = 12.2
$ws.Range("E46").Value = "no"

$loInvoices = $ws.ListObjects.Add(1, $ws.Range("A41:F49"), 0, 1)
$loInvoices.Name = "Table1345"
$loInvoices.TableStyle = "TableStyleMedium14"

# ---------------------------------------------------------------------------
# Final selection to reflect where the author ended up editing
# ---------------------------------------------------------------------------
$ws.Range("F48").Select()
